$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1209.0999999999999
$ws.Range("H2").Value = 2

$ws.Range("H7").Select()
